# Populate Sheet1 with the Name/Age/City/Address table that was blank before.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(1, "John", 25, "New York", "25th street"),
    @(2, "Alice", 30, "Los Angeles", "apartment 25, 5th floor"),
    @(3, "Bob", 22, "Chicago", "26th street"),
    @(4, "Charlie", 28, "Houston", "apartment 25, 5th floor"),
    @(5, "David", 35, "Phoenix", "27th street"),
    @(6, "Emma", 40, "Philadelphia", "apartment 25, 5th floor"),
    @(7, "Fathi", 19, "San Antonio", "28th street"),
    @(8, "Grace", 21, "San Diego", "apartment 25, 5th floor"),
    @(9, "Henry", 45, "Dallas", "29th street"),
    @(10, "Ivy", 50, "San Jose", "apartment 25, 5th floor"),
    @(11, "Jack", 33, "Austin", "30th street"),
    @(12, "Karen", 27, "Jacksonville", "apartment 25, 5th floor"),
    @(13, "Leo", 31, "Fort Worth", "31st street"),
    @(14, "Mona", 29, "Columbus", "apartment 25, 5th floor"),
    @(15, "Nathan", 26, "Charlotte", "32nd street"),
    @(16, "Olivia", 38, "Indianapolis", "apartment 25, 5th floor"),
    @(17, "Paul", 41, "Seattle", "33rd street"),
    @(18, "Quinn", 20, "Denver", "apartment 25, 5th floor"),
    @(19, "Rachel", 32, "Washington", "34th street"),
    @(20, "Steve", 37, "Boston", "apartment 25, 5th floor")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$wb.Save()
